$d = $word.ActiveDocument

$d.Content.Find.Execute("509÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "628÷2=", 2) | Out-Null
$d.Content.Find.Execute("554÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "302÷9=", 2) | Out-Null
$d.Content.Find.Execute("586÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "780÷4=", 2) | Out-Null
$d.Content.Find.Execute("447÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "727÷9=", 2) | Out-Null
$d.Content.Find.Execute("712÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "197÷6=", 2) | Out-Null
$d.Content.Find.Execute("379÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "394÷9=", 2) | Out-Null
$d.Content.Find.Execute("824÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "177÷3=", 2) | Out-Null
$d.Content.Find.Execute("731÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "931÷5=", 2) | Out-Null
$d.Content.Find.Execute("315÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "509÷7=", 2) | Out-Null
$d.Content.Find.Execute("488÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "398÷3=", 2) | Out-Null
$d.Content.Find.Execute("682÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "566÷4=", 2) | Out-Null
$d.Content.Find.Execute("691÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "114÷8=", 2) | Out-Null
$d.Content.Find.Execute("944÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "990÷4=", 2) | Out-Null
$d.Content.Find.Execute("121÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "757÷4=", 2) | Out-Null
$d.Content.Find.Execute("552÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "407÷3=", 2) | Out-Null
$d.Content.Find.Execute("453÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "521÷4=", 2) | Out-Null
$d.Content.Find.Execute("585÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "477÷3=", 2) | Out-Null
$d.Content.Find.Execute("858÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "514÷4=", 2) | Out-Null
$d.Content.Find.Execute("547÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "510÷3=", 2) | Out-Null
$d.Content.Find.Execute("672÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "173÷6=", 2) | Out-Null
$d.Content.Find.Execute("324÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "755÷3=", 2) | Out-Null
$d.Content.Find.Execute("340÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "230÷4=", 2) | Out-Null
$d.Content.Find.Execute("780÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "648÷7=", 2) | Out-Null
$d.Content.Find.Execute("611÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "137÷4=", 2) | Out-Null
$d.Content.Find.Execute("629÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "270÷7=", 2) | Out-Null
